$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.103.16"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").Value = "1.640.87"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.57%  "

$ws.Range("E6").Value = "  +1.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.33"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.80%  "

$ws.Range("E9").Value = "  +4.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0615"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0916"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "

$ws.Range("D12").Value = "1.878.60"
$ws.Range("E12").Value = "  +2.66%  "

$ws.Range("D13").Value = "1.645.17"
$ws.Range("E13").Value = "  +2.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.576"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +24.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.91"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.25%  "

$ws.Range("D17").Value = "30.126.25"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.95%  "

$ws.Range("D20").Value = "0.0₃0710"
$ws.Range("E20").Value = "  +2.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.37%  "

$ws.Range("E24").Value = "  +1.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.42%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.111"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.03%  "

$ws.Range("E28").Value = "  +4.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0493"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.84%  "

$ws.Range("E31").Value = "  +6.20%  "

$ws.Range("E32").Value = "  +5.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").Value = "1.438.19"
$ws.Range("E34").Value = "  +1.24%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.92%  "

$ws.Range("E36").Value = "  +1.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("B38").Value = "Aave"
$ws.Range("C38").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "77.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +17.65%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0172"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.562"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.03%  "

$ws.Range("E42").Value = "  +3.10%  "

$ws.Range("E43").Value = "  +3.30%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.78%  "

$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("E46").Value = "  +5.35%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "

$ws.Range("D49").Value = "1.785.47"
$ws.Range("E49").Value = "  +2.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.79%  "

$ws.Range("D51").Value = "0.0₆0112"
$ws.Range("E51").Value = "  +7.65%  "

